$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.605.28'
$ws.Range("E2").Value = '  -2.55%  '
$ws.Range("D3").Value = '1.981.75'
$ws.Range("E3").Value = '  -3.64%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.41'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.636'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.83'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +6.02%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '58.47'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.360'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0734'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("E12").Value = '  -2.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.947'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("E14").Value = '  -1.52%  '
$ws.Range("D15").Value = '2.268.05'
$ws.Range("E15").Value = '  -3.77%  '
$ws.Range("E16").Value = '  -2.33%  '
$ws.Range("D17").Value = '1.973.06'
$ws.Range("E17").Value = '  -4.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.01'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +7.00%  '
$ws.Range("D19").Value = '35.565.09'
$ws.Range("E19").Value = '  -2.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.37'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").Value = '0.0₃0842'
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '233.02'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.17'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.55'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +20.09%  '
$ws.Range("E26").Value = '  -2.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.65'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.09'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.17'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.51%  '
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.86'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.11'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0951'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +15.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0592'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("E35").Value = '  +9.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.33'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.44%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.76'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.22'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +7.81%  '
$ws.Range("E40").Value = '  -2.25%  '
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0212'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.80%  '
$ws.Range("E43").Value = '  -1.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '91.66'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.54%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.60'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("E47").Value = '  -4.78%  '
$ws.Range("D48").Value = '1.369.00'
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.89'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.37%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '47.55'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.03%  '
